# Fruta / hortaliza, semanal
# Insert a new weekly data row at row 27 (pushing the existing rows 27..129
# down to 28..130) and populate it with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 27; Excel shifts rows 27-129 down
# to 28-130 and the sheet's used range grows to R130 automatically.
$ws.Rows.Item(27).Insert()

# Populate the newly inserted row 27 with the new weekly record.
$ws.Range("A27").Value2 = 5
$ws.Range("B27").Value2 = "Macroferia Regional de Talca"
$ws.Range("C27").Value2 = "Maule"
$ws.Range("D27").Value2 = 44910
$ws.Range("E27").Value2 = 7
$ws.Range("F27").Value2 = 100112022
$ws.Range("G27").Value2 = "Arveja Verde"
$ws.Range("H27").Value2 = "Sin especificar"
$ws.Range("I27").Value2 = "Primera"
$ws.Range("J27").Value2 = 500
$ws.Range("K27").Value2 = 22000
$ws.Range("L27").Value2 = 22000
$ws.Range("M27").Value2 = 22000
$ws.Range("N27").Value2 = "`$/saco 25 kilos"
$ws.Range("O27").Value2 = "Carahue"
$ws.Range("P27").Value2 = 880
$ws.Range("Q27").Value2 = 25
$ws.Range("R27").Value2 = "Hortaliza"
